$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 715, shifting existing rows 715:736 down to 716:737
$ws.Rows.Item(715).Insert()

# Populate the newly inserted row 715 with data
$ws.Cells.Item(715, 1).Value = 10
$ws.Cells.Item(715, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(715, 3).Value = "La Araucanía"
$ws.Cells.Item(715, 4).Value = 45075
$ws.Cells.Item(715, 5).Value = 9
$ws.Cells.Item(715, 6).Value = 100112028
$ws.Cells.Item(715, 7).Value = "Sandia"
$ws.Cells.Item(715, 8).Value = "Sin especificar"
$ws.Cells.Item(715, 9).Value = "Primera"
$ws.Cells.Item(715, 10).Value = 500
$ws.Cells.Item(715, 11).Value = 3500
$ws.Cells.Item(715, 12).Value = 3500
$ws.Cells.Item(715, 13).Value = 3500
$ws.Cells.Item(715, 14).Value = "`$/unidad"
$ws.Cells.Item(715, 15).Value = "Brasil"
$ws.Cells.Item(715, 16).Value = 3500
$ws.Cells.Item(715, 17).Value = 1
$ws.Cells.Item(715, 18).Value = "Hortaliza"

# Copy the date cell style (s="2") from the row below (old row 715, now at 716) to preserve formatting
$ws.Cells.Item(716, 4).Copy()
$ws.Cells.Item(715, 4).PasteSpecial(-4122)
$ws.Cells.Item(715, 4).Value = 45075
